# Add calculation of capacity in MWe
# Adds a new "Size (MWe)" column (D) to the ReverseEngineer sheet, computed
# from an electricity-consumption rate (B7, MWhe/m3-H2O) times the plant
# size (E column, m3/h), for both cost blocks (rows 7-16 and rows 43-52).
# The downstream log-log regression cells (E14:E16 and E50:E52) are
# re-pointed to use the new MWe-based size instead of the raw m3/h size.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ReverseEngineer")

# ---- Block 1 header (row 7) ----------------------------------------------
# New header cells: D7 first (so it becomes shared-string index 51), then
# A7, then C7, matching the ordering Excel used when it originally saved
# these new labels.
$ws.Range("D7").Value = "Size (MWe)"
$ws.Range("A7").Value = "Electricity consumption"
$ws.Range("C7").Value = "MWhe/m3-H2O"

# Electricity-consumption rate used by both blocks (MWhe per m3 of water).
$ws.Range("B7").Value = 0.00325

# ---- Block 1 data (rows 8-11): Size (MWe) = $B$7 * E -----------------------
$ws.Range("D8").NumberFormat = "0.000"
$ws.Range("D8").Formula = '=$B$7*E8'

$ws.Range("D9:D11").NumberFormat = "0.000"
$ws.Range("D9").Formula = '=$B$7*E9'
$ws.Range("D10").Formula = '=$B$7*E10'
$ws.Range("D11").Formula = '=$B$7*E11'

# ---- Block 1 regression (rows 14-16) now reference column D ---------------
$ws.Range("E14").Formula = '=LOG(D9/$D$9)'
$ws.Range("E15").Formula = '=LOG(D10/$D$9)'
$ws.Range("E16").Formula = '=LOG(D11/$D$9)'

# ---- Block 2 header (row 43) ----------------------------------------------
$ws.Range("D43").Value = "Size (MWe)"

# ---- Block 2 data (rows 44-47): Size (MWe) = E * $B$7 ----------------------
$ws.Range("D44").NumberFormat = "0.000"
$ws.Range("D44").Formula = '=E44*$B$7'

$ws.Range("D45:D47").NumberFormat = "0.000"
$ws.Range("D45").Formula = '=E45*$B$7'
$ws.Range("D46").Formula = '=E46*$B$7'
$ws.Range("D47").Formula = '=E47*$B$7'

# ---- Block 2 regression (rows 50-52) now reference column D ---------------
$ws.Range("E50").Formula = '=LOG(D45/$D$45)'
$ws.Range("E51").Formula = '=LOG(D46/$D$45)'
$ws.Range("E52").Formula = '=LOG(D47/$D$45)'

# ---- Reposition/resize the second (CAPEX) chart ---------------------------
# It moved up and slightly right/wider, from anchoring near row 39 to
# anchoring near row 28, to make room below it.
$co = $ws.ChartObjects().Item(2)
$co.Left = 581.8477349901575
$co.Top = 429.1429133858268
$co.Width = 384.88599870816927
$co.Height = 216.0

# ---- Restore the view: selection on G21, scrolled to top-left -------------
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("G21").Select()

$wb.Save()
